# Aerodynamic database update: the Xcg ESTIMATION METHOD COMPARISON
# tables swap the relative order of the "TORENBEEK_1982" and "SFORZA"
# rows (and their associated values) on the FUSELAGE and WING sheets.
#
# FUSELAGE: single comparison table (rows 11-12)
# WING:     two comparison tables (rows 11-12 and rows 15-16)

$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet ---------------------------------------------------
$fuselage = $wb.Worksheets.Item("FUSELAGE")

$fuselage.Range("A11").Value = "SFORZA"
$fuselage.Range("C11").Value = 17.143322222222217
$fuselage.Range("A12").Value = "TORENBEEK_1982"
$fuselage.Range("C12").Value = 16.8345

# --- WING sheet ---------------------------------------------------------
$wing = $wb.Worksheets.Item("WING")

# Xcg ESTIMATION METHOD COMPARISON table
$wing.Range("A11").Value = "SFORZA"
$wing.Range("C11").Value = 4.3629715646212155
$wing.Range("A12").Value = "TORENBEEK_1982"
$wing.Range("C12").Value = 3.5180298935880643

# Ycg ESTIMATION METHOD COMPARISON table
$wing.Range("A15").Value = "SFORZA"
$wing.Range("C15").Value = 4.998846772296348
$wing.Range("A16").Value = "TORENBEEK_1982"
$wing.Range("C16").Value = 6.114221148470394
